$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "MISSING: $find"
    }
}

# --- Paragraph: genetic algorithm description ---
Replace-Text "mimicking the evolution of life in nature on a suitably simplistic scale" "mimicking the natural selection on a suitably simplistic scale"
Replace-Text "are selected for breeding to create the next generation." "are selected for creating the next generation."
Replace-Text "Additional factors like mutation in chromosomes and crossovers among subpopulations may also be specified to reduce" "Additional factors like random mutation in genomes are also specified to reduce"

# --- Paragraph: deployment strategy ---
Replace-Text "the node arrangement is allowed to expand in the desired aspect ratio." "the node arrangement is allowed to expand within these bounds."
Replace-Text "and at best might possibly contribute to face coverage and alternate communication routes. In the case" "and at best contribute to face coverage and alternate communication routes. A suitable fraction of the score is attributed to every node lying on the convex hull, and is used as the unit penalty score for every node that violates the specified bounds. In using this approach, nodes that lie within the convex hull, and thus do not contribute to the score of the individual, increase the unit bounds violation penalty, discouraging such arrangements as a side effect. In the case"

# --- Simplify empty rPr in math paragraphs (sz/szCs 24 -> empty) ---
# Handled via direct run property clears below.

# --- New character styles ---
$s1 = $d.Styles.Add("ListLabel 26", 2)
$s1.Font.NameBi = "Symbol"
$s2 = $d.Styles.Add("ListLabel 27", 2)
$s2.Font.Size = 10
$s2.Font.SizeBi = 10
$s3 = $d.Styles.Add("ListLabel 28", 2)
$s3.Font.NameBi = "OpenSymbol"
